$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp shown at the top of the sheet
$ws.Range("A1").Value = "Datos actualizados a 14 de Agosto de 2020 a las 16:40"

# Refresh per-country case counts; a handful of countries changed rank
# (new totals re-sort the table), so both the country name and the
# numeric columns are rewritten for the affected rows.

# Row 4: Estados Unidos
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 5418071
$ws.Range("C4").Value = 2405
$ws.Range("D4").Value = 2844262
$ws.Range("E4").Value = 2403363
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 31
$ws.Range("H4").Value = 170446

# Row 6: India
$ws.Range("A6").Value = "India"
$ws.Range("B6").Value = 2483567
$ws.Range("C6").Value = 23954
$ws.Range("D6").Value = 1770682
$ws.Range("E6").Value = 664518
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 223
$ws.Range("H6").Value = 48367

# Row 18: Argentina
$ws.Range("A18").Value = "Argentina"
$ws.Range("B18").Value = 276072
$ws.Range("C18").Value = 0
$ws.Range("D18").Value = 199005
$ws.Range("E18").Value = 71639
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 66
$ws.Range("H18").Value = 5428

# Row 22: Alemania
$ws.Range("A22").Value = "Alemania"
$ws.Range("B22").Value = 222487
$ws.Range("C22").Value = 218
$ws.Range("D22").Value = 200800
$ws.Range("E22").Value = 12404
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 9283

# Row 36: Suecia
$ws.Range("A36").Value = "Suecia"
$ws.Range("B36").Value = 84294
$ws.Range("C36").Value = 0
$ws.Range("D36").Value = 0
$ws.Range("E36").Value = 0
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 7
$ws.Range("H36").Value = 5783

# Row 61: Azerbaiyan
$ws.Range("A61").Value = "Azerbaiyan"
$ws.Range("B61").Value = 34018
$ws.Range("C61").Value = 103
$ws.Range("D61").Value = 31490
$ws.Range("E61").Value = 2024
$ws.Range("F61").Value = 0
$ws.Range("G61").Value = 4
$ws.Range("H61").Value = 504

# Row 64: Moldavia
$ws.Range("A64").Value = "Moldavia"
$ws.Range("B64").Value = 29483
$ws.Range("C64").Value = 396
$ws.Range("D64").Value = 20556
$ws.Range("E64").Value = 8043
$ws.Range("F64").Value = 0
$ws.Range("G64").Value = 6
$ws.Range("H64").Value = 884

# Row 65: Kenia
$ws.Range("A65").Value = "Kenia"
$ws.Range("B65").Value = 29334
$ws.Range("C65").Value = 580
$ws.Range("D65").Value = 15298
$ws.Range("E65").Value = 13571
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 5
$ws.Range("H65").Value = 465

# Row 66: Serbia
$ws.Range("A66").Value = "Serbia"
$ws.Range("B66").Value = 29233
$ws.Range("C66").Value = 235
$ws.Range("D66").Value = 26117
$ws.Range("E66").Value = 2451
$ws.Range("F66").Value = 0
$ws.Range("G66").Value = 4
$ws.Range("H66").Value = 665

# Row 83: Republica de Macedonia
$ws.Range("A83").Value = "Republica de Macedonia"
$ws.Range("B83").Value = 12515
$ws.Range("C83").Value = 158
$ws.Range("D83").Value = 9030
$ws.Range("E83").Value = 2950
$ws.Range("F83").Value = 0
$ws.Range("G83").Value = 3
$ws.Range("H83").Value = 535

# Row 86: Noruega
$ws.Range("A86").Value = "Noruega"
$ws.Range("B86").Value = 9887
$ws.Range("C86").Value = 36
$ws.Range("D86").Value = 8857
$ws.Range("E86").Value = 769
$ws.Range("F86").Value = 0
$ws.Range("G86").Value = 4
$ws.Range("H86").Value = 261

# Row 143: Uganda
$ws.Range("A143").Value = "Uganda"
$ws.Range("B143").Value = 1385
$ws.Range("C143").Value = 32
$ws.Range("D143").Value = 1142
$ws.Range("E143").Value = 231
$ws.Range("F143").Value = 0
$ws.Range("G143").Value = 1
$ws.Range("H143").Value = 12

# Row 174: Belice
$ws.Range("A174").Value = "Belice"
$ws.Range("B174").Value = 356
$ws.Range("C174").Value = 60
$ws.Range("D174").Value = 32
$ws.Range("E174").Value = 322
$ws.Range("F174").Value = 0
$ws.Range("G174").Value = 0
$ws.Range("H174").Value = 2

# Row 175: Mauricio
$ws.Range("A175").Value = "Mauricio"
$ws.Range("B175").Value = 344
$ws.Range("C175").Value = 0
$ws.Range("D175").Value = 334
$ws.Range("E175").Value = 0
$ws.Range("F175").Value = 0
$ws.Range("G175").Value = 0
$ws.Range("H175").Value = 10

# Row 176: Martinica
$ws.Range("A176").Value = "Martinica"
$ws.Range("B176").Value = 336
$ws.Range("C176").Value = 0
$ws.Range("D176").Value = 98
$ws.Range("E176").Value = 222
$ws.Range("F176").Value = 0
$ws.Range("G176").Value = 0
$ws.Range("H176").Value = 16

# Row 177: Isla de Man
$ws.Range("A177").Value = "Isla de Man"
$ws.Range("B177").Value = 336
$ws.Range("C177").Value = 0
$ws.Range("D177").Value = 312
$ws.Range("E177").Value = 0
$ws.Range("F177").Value = 0
$ws.Range("G177").Value = 0
$ws.Range("H177").Value = 24

# Row 178: Mongolia
$ws.Range("A178").Value = "Mongolia"
$ws.Range("B178").Value = 297
$ws.Range("C178").Value = 0
$ws.Range("D178").Value = 269
$ws.Range("E178").Value = 28
$ws.Range("F178").Value = 0
$ws.Range("G178").Value = 0
$ws.Range("H178").Value = 0

# Row 195: Liechtenstein
$ws.Range("A195").Value = "Liechtenstein"
$ws.Range("B195").Value = 91
$ws.Range("C195").Value = 1
$ws.Range("D195").Value = 87
$ws.Range("E195").Value = 3
$ws.Range("F195").Value = 0
$ws.Range("G195").Value = 0
$ws.Range("H195").Value = 1

# Row 213: Islas Malvinas
$ws.Range("A213").Value = "Islas Malvinas"
$ws.Range("B213").Value = 13
$ws.Range("C213").Value = 0
$ws.Range("D213").Value = 13
$ws.Range("E213").Value = 0
$ws.Range("F213").Value = 0
$ws.Range("G213").Value = 0
$ws.Range("H213").Value = 0

# Row 214: Montserrat
$ws.Range("A214").Value = "Montserrat"
$ws.Range("B214").Value = 13
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 12
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 1
